{"js": "// The document holds one \"header\" paragraph (date line) followed by a\n// table. Most table cells are blank; the populated cells (and the date\n// paragraph) carry the text that needs updating. Several of the source\n// strings repeat verbatim but map to DIFFERENT replacement strings\n// depending on their position, so we can't do a blind global\n// find/replace - we walk every paragraph in document order and replace\n// by position among the non-empty ones.\n\nconst replacements = [\n  \"2024-06-22 Saturday\",\n  \"41\u00f74=10, 1\",\n  \"45\u00f76=7, 3\",\n  \"74\u00f79=8, 2\",\n  \"47\u00f79=5, 2\",\n  \"69\u00f74=17, 1\",\n  \"27\u00f76=4, 3\",\n  \"74\u00f78=9, 2\",\n  \"77\u00f76=12, 5\",\n  \"64\u00f73=21, 1\",\n  \"56\u00f74=14, 0\",\n  \"42\u00f72=21, 0\",\n  \"98\u00f77=14, 0\",\n  \"75\u00f79=8, 3\",\n  \"26\u00f75=5, 1\",\n  \"51\u00f79=5, 6\",\n  \"70\u00f76=11, 4\",\n  \"51\u00f79=5, 6\",\n  \"29\u00f79=3, 2\",\n  \"73\u00f77=10, 3\",\n  \"34\u00f72=17, 0\",\n  \"13\u00f74=3, 1\",\n  \"64\u00f74=16, 0\",\n  \"29\u00f78=3, 5\",\n  \"93\u00f76=15, 3\",\n  \"50\u00f73=16, 2\",\n];\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet idx = 0;\nfor (const para of paragraphs.items) {\n  if (para.text !== \"\") {\n    if (idx >= replacements.length) {\n      throw new Error(\"More populated paragraphs than replacements available\");\n    }\n    para.insertText(replacements[idx], Word.InsertLocation.replace);\n    idx++;\n  }\n}\n\nif (idx !== replacements.length) {\n  throw new Error(`Expected ${replacements.length} populated paragraphs, found ${idx}`);\n}\n\nawait context.sync();\n", "ps1": "# The document holds one \"header\" paragraph (date line) followed by a\n# single 20x5 table. Only every 5th row (rows 1, 5, 9, 13, 17) carries\n# text; the rest of the cells are blank. Several of the source cell\n# strings repeat verbatim but map to DIFFERENT replacement strings\n# depending on their position (e.g. \"46\u00f78=5, 6\" appears twice and\n# becomes \"45\u00f76=7, 3\" the first time and \"51\u00f79=5, 6\" the second), so a\n# blind document-wide Find/Replace would be wrong. Instead we replace\n# the date line directly, then walk the table's cells in row-major\n# order and replace each non-blank one by position.\n\n$d = $word.ActiveDocument\n\n# --- date paragraph -------------------------------------------------\n# Trim the trailing paragraph mark first so the assignment below only\n# replaces the visible text (assigning text with an embedded `r would\n# split this into two paragraphs instead of renaming the one line).\n$dateRange = $d.Paragraphs.Item(1).Range\n$dateRange.MoveEnd(1, -1) | Out-Null\n$dateRange.Text = \"2024-06-22 Saturday\"\n\n# --- table cells ------------------------------------------------------\n$cellReplacements = @(\n    \"41\u00f74=10, 1\", \"45\u00f76=7, 3\", \"74\u00f79=8, 2\", \"47\u00f79=5, 2\", \"69\u00f74=17, 1\",\n    \"27\u00f76=4, 3\", \"74\u00f78=9, 2\", \"77\u00f76=12, 5\", \"64\u00f73=21, 1\", \"56\u00f74=14, 0\",\n    \"42\u00f72=21, 0\", \"98\u00f77=14, 0\", \"75\u00f79=8, 3\", \"26\u00f75=5, 1\", \"51\u00f79=5, 6\",\n    \"70\u00f76=11, 4\", \"51\u00f79=5, 6\", \"29\u00f79=3, 2\", \"73\u00f77=10, 3\", \"34\u00f72=17, 0\",\n    \"13\u00f74=3, 1\", \"64\u00f74=16, 0\", \"29\u00f78=3, 5\", \"93\u00f76=15, 3\", \"50\u00f73=16, 2\"\n)\n\n$t = $d.Tables.Item(1)\n$idx = 0\nfor ($row = 1; $row -le $t.Rows.Count; $row++) {\n    for ($col = 1; $col -le $t.Columns.Count; $col++) {\n        $cell = $t.Cell($row, $col)\n        $r = $cell.Range\n        # Strip the trailing cell-mark/paragraph-mark pair so the\n        # assignment below only touches the visible text. Use Start/End\n        # (not .Text, which can read stale/misleading content on a\n        # collapsed range) to detect whether anything real is left.\n        $r.MoveEnd(1, -1) | Out-Null\n        if ($r.End -gt $r.Start) {\n            if ($idx -ge $cellReplacements.Count) {\n                throw \"More populated cells than replacements available\"\n            }\n            $r.Text = $cellReplacements[$idx]\n            $idx++\n        }\n    }\n}\n\nif ($idx -ne $cellReplacements.Count) {\n    throw \"Expected $($cellReplacements.Count) populated cells, found $idx\"\n}\n"}
